$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Generate Report for Handback
#
# The localization-status report is regenerated: the three tracked source
# files (ffff7664ed3f..., ffffff470004c4..., cad7c9f4...) have all reached
# "Handed back: in sync with en-US" status, and a fresh handback round-trip
# for cad7c9f4 (handoff 03:26:49 -> handback 03:26:56/03:27:16 for zh-cn,
# 03:27:24 for de-de) has been folded into the report rows.
# ---------------------------------------------------------------------------

$statusHandedBack = "Handed back: in sync with en-US"
$include = "Include"
$mdExt = ".md"

# -- URLs reused verbatim from the original workbook (by hyperlink slot) ----
$url_md_ffff7664 = "https://github.com/OpenLocalizationTest/oltest/blob/2cb5198c13a4dc417c4e417b71266e29c9d9992c/e2e/ffff7664ed3f-a6c2-434a-9dd7-0b212984201f.md"
$url_md_ffffff47 = "https://github.com/OpenLocalizationTest/oltest/blob/2cb5198c13a4dc417c4e417b71266e29c9d9992c/e2e/ffffff470004c4-cb39-4dcf-8882-498e14e937c9.md"
$url_md_cad7c9f4 = "https://github.com/OpenLocalizationTest/oltest/blob/2cb5198c13a4dc417c4e417b71266e29c9d9992c/e2e/cad7c9f4-99db-479c-9fce-e95c53336e9e.md"

$url_zhcn_handoff_b0489487 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1e9b2c09f568629c5009f45ddfea0df21ad0deac/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b0489487-19e5-4b3f-87d8-aaa577d701b1.76e6209aa82c7c87aa7149c3caddf6768b87afab.zh-cn.xlf"
$url_zhcn_target_b0489487  = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/aac5fcc52e2b3f5ac80c1059a5b055fc87b7ad68/e2e/b0489487-19e5-4b3f-87d8-aaa577d701b1.md"
$url_zhcn_handback_b0489487 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/bea60ad2571e96987ca596a3279d8bf314985384/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b0489487-19e5-4b3f-87d8-aaa577d701b1.76e6209aa82c7c87aa7149c3caddf6768b87afab.zh-cn.xlf"
$url_zhcn_handoff_cad7c9f4 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/67da43f9a9e63651cdabacd504ef5f0c78d81cf3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/cad7c9f4-99db-479c-9fce-e95c53336e9e.e92a848c82665bf95872cb8eebbcd5ad40d40087.zh-cn.xlf"

$url_dede_handoff_b0489487 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c9b03d97ced3134a59132f504c427f1c0670f0ae/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b0489487-19e5-4b3f-87d8-aaa577d701b1.76e6209aa82c7c87aa7149c3caddf6768b87afab.de-de.xlf"
$url_dede_target_b0489487  = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/5f0438e4bb2082398d1db1e29524edc6dfe80afc/e2e/b0489487-19e5-4b3f-87d8-aaa577d701b1.md"
$url_dede_handback_b0489487 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/d2ac6cb0aad422321972fe201be986b57518895e/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b0489487-19e5-4b3f-87d8-aaa577d701b1.76e6209aa82c7c87aa7149c3caddf6768b87afab.de-de.xlf"
$url_dede_handoff_cad7c9f4 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/460e34a9bf60d59abb47a99ac8da8116a66384a8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/cad7c9f4-99db-479c-9fce-e95c53336e9e.e92a848c82665bf95872cb8eebbcd5ad40d40087.de-de.xlf"

# ===========================================================================
# Sheet "Overview"
# ===========================================================================
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("B2").Value = $statusHandedBack
$ov.Range("C2").Value = $statusHandedBack
$ov.Range("D2").Value = "2016-03-25 03:26:56"

$ov.Range("B3").Value = $statusHandedBack
$ov.Range("C3").Value = $statusHandedBack
$ov.Range("D3").Value = "2016-03-25 03:24:57"

$ov.Range("B4").Value = $statusHandedBack
$ov.Range("C4").Value = $statusHandedBack
$ov.Range("D4").Value = "2016-03-25 03:24:57"

$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), $url_md_ffff7664, [Type]::Missing, [Type]::Missing, "cad7c9f4-99db-479c-9fce-e95c53336e9e.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A3"), $url_md_ffffff47, [Type]::Missing, [Type]::Missing, "ffff7664ed3f-a6c2-434a-9dd7-0b212984201f.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A4"), $url_md_cad7c9f4, [Type]::Missing, [Type]::Missing, "ffffff470004c4-cb39-4dcf-8882-498e14e937c9.md") | Out-Null

# ===========================================================================
# Sheet "zh-cn"
# ===========================================================================
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("B2").Value = $mdExt
$zh.Range("C2").Value = $statusHandedBack
$zh.Range("E2").Value = "2016-03-25 03:26:49"
$zh.Range("H2").Value = "2016-03-25 03:27:16"
$zh.Range("J2").Value = $include

$zh.Range("B3").Value = $mdExt
$zh.Range("C3").Value = $statusHandedBack
$zh.Range("E3").Value = "2016-03-25 03:24:53"
$zh.Range("H3").Value = "2016-03-25 03:25:19"
$zh.Range("J3").Value = $include

$zh.Range("A4").Value = "ffffff470004c4-cb39-4dcf-8882-498e14e937c9.md"
$zh.Range("B4").Value = $mdExt
$zh.Range("C4").Value = $statusHandedBack
$zh.Range("E4").Value = "2016-03-25 03:24:53"
$zh.Range("H4").Value = "2016-03-25 03:25:19"
$zh.Range("J4").Value = $include

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), $url_md_ffff7664, [Type]::Missing, [Type]::Missing, "cad7c9f4-99db-479c-9fce-e95c53336e9e.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D2"), $url_zhcn_handoff_b0489487, [Type]::Missing, [Type]::Missing, "cad7c9f4-99db-479c-9fce-e95c53336e9e.e92a848c82665bf95872cb8eebbcd5ad40d40087.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F2"), $url_zhcn_target_b0489487, [Type]::Missing, [Type]::Missing, "cad7c9f4-99db-479c-9fce-e95c53336e9e.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("G2"), $url_zhcn_handback_b0489487, [Type]::Missing, [Type]::Missing, "cad7c9f4-99db-479c-9fce-e95c53336e9e.e92a848c82665bf95872cb8eebbcd5ad40d40087.zh-cn.xlf") | Out-Null

$zh.Hyperlinks.Add($zh.Range("A3"), $url_md_ffffff47, [Type]::Missing, [Type]::Missing, "ffff7664ed3f-a6c2-434a-9dd7-0b212984201f.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D3"), $url_zhcn_handoff_b0489487, [Type]::Missing, [Type]::Missing, "b0489487-19e5-4b3f-87d8-aaa577d701b1.76e6209aa82c7c87aa7149c3caddf6768b87afab.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F3"), $url_zhcn_target_b0489487, [Type]::Missing, [Type]::Missing, "b0489487-19e5-4b3f-87d8-aaa577d701b1.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("G3"), $url_zhcn_handback_b0489487, [Type]::Missing, [Type]::Missing, "b0489487-19e5-4b3f-87d8-aaa577d701b1.76e6209aa82c7c87aa7149c3caddf6768b87afab.zh-cn.xlf") | Out-Null

$zh.Hyperlinks.Add($zh.Range("A4"), $url_md_cad7c9f4, [Type]::Missing, [Type]::Missing, "ffffff470004c4-cb39-4dcf-8882-498e14e937c9.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D4"), $url_zhcn_handoff_b0489487, [Type]::Missing, [Type]::Missing, "b0489487-19e5-4b3f-87d8-aaa577d701b1.76e6209aa82c7c87aa7149c3caddf6768b87afab.zh-cn.xlf") | Out-Null
$zh.Hyperlinks.Add($zh.Range("F4"), $url_zhcn_target_b0489487, [Type]::Missing, [Type]::Missing, "b0489487-19e5-4b3f-87d8-aaa577d701b1.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("G4"), $url_zhcn_handback_b0489487, [Type]::Missing, [Type]::Missing, "b0489487-19e5-4b3f-87d8-aaa577d701b1.76e6209aa82c7c87aa7149c3caddf6768b87afab.zh-cn.xlf") | Out-Null

# ===========================================================================
# Sheet "de-de"
# ===========================================================================
$de = $wb.Worksheets.Item("de-de")

$de.Range("B2").Value = $mdExt
$de.Range("C2").Value = $statusHandedBack
$de.Range("E2").Value = "2016-03-25 03:26:56"
$de.Range("H2").Value = "2016-03-25 03:27:24"
$de.Range("J2").Value = $include

$de.Range("B3").Value = $mdExt
$de.Range("C3").Value = $statusHandedBack
$de.Range("E3").Value = "2016-03-25 03:24:57"
$de.Range("H3").Value = "2016-03-25 03:25:27"
$de.Range("J3").Value = $include

$de.Range("A4").Value = "ffffff470004c4-cb39-4dcf-8882-498e14e937c9.md"
$de.Range("B4").Value = $mdExt
$de.Range("C4").Value = $statusHandedBack
$de.Range("E4").Value = "2016-03-25 03:24:57"
$de.Range("H4").Value = "2016-03-25 03:25:27"
$de.Range("J4").Value = $include

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), $url_md_ffff7664, [Type]::Missing, [Type]::Missing, "cad7c9f4-99db-479c-9fce-e95c53336e9e.md") | Out-Null
$de.Hyperlinks.Add($de.Range("D2"), $url_dede_handoff_b0489487, [Type]::Missing, [Type]::Missing, "cad7c9f4-99db-479c-9fce-e95c53336e9e.e92a848c82665bf95872cb8eebbcd5ad40d40087.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("F2"), $url_dede_target_b0489487, [Type]::Missing, [Type]::Missing, "cad7c9f4-99db-479c-9fce-e95c53336e9e.md") | Out-Null
$de.Hyperlinks.Add($de.Range("G2"), $url_dede_handback_b0489487, [Type]::Missing, [Type]::Missing, "cad7c9f4-99db-479c-9fce-e95c53336e9e.e92a848c82665bf95872cb8eebbcd5ad40d40087.de-de.xlf") | Out-Null

$de.Hyperlinks.Add($de.Range("A3"), $url_md_ffffff47, [Type]::Missing, [Type]::Missing, "ffff7664ed3f-a6c2-434a-9dd7-0b212984201f.md") | Out-Null
$de.Hyperlinks.Add($de.Range("D3"), $url_dede_handoff_b0489487, [Type]::Missing, [Type]::Missing, "b0489487-19e5-4b3f-87d8-aaa577d701b1.76e6209aa82c7c87aa7149c3caddf6768b87afab.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("F3"), $url_dede_target_b0489487, [Type]::Missing, [Type]::Missing, "b0489487-19e5-4b3f-87d8-aaa577d701b1.md") | Out-Null
$de.Hyperlinks.Add($de.Range("G3"), $url_dede_handback_b0489487, [Type]::Missing, [Type]::Missing, "b0489487-19e5-4b3f-87d8-aaa577d701b1.76e6209aa82c7c87aa7149c3caddf6768b87afab.de-de.xlf") | Out-Null

$de.Hyperlinks.Add($de.Range("A4"), $url_md_cad7c9f4, [Type]::Missing, [Type]::Missing, "ffffff470004c4-cb39-4dcf-8882-498e14e937c9.md") | Out-Null
$de.Hyperlinks.Add($de.Range("D4"), $url_dede_handoff_b0489487, [Type]::Missing, [Type]::Missing, "b0489487-19e5-4b3f-87d8-aaa577d701b1.76e6209aa82c7c87aa7149c3caddf6768b87afab.de-de.xlf") | Out-Null
$de.Hyperlinks.Add($de.Range("F4"), $url_dede_target_b0489487, [Type]::Missing, [Type]::Missing, "b0489487-19e5-4b3f-87d8-aaa577d701b1.md") | Out-Null
$de.Hyperlinks.Add($de.Range("G4"), $url_dede_handback_b0489487, [Type]::Missing, [Type]::Missing, "b0489487-19e5-4b3f-87d8-aaa577d701b1.76e6209aa82c7c87aa7149c3caddf6768b87afab.de-de.xlf") | Out-Null

$ov.Select()
